$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("B1").Value = "http://dbpedia.org/ontology/deathPlace"
$ws.Range("C1").Value = "http://dbpedia.org/ontology/parent"
$ws.Range("D1").Value = "http://dbpedia.org/ontology/deathDate"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/birthDate"
$ws.Range("F1").Value = "http://dbpedia.org/ontology/birthPlace"

# Carry the header style (bold/border/centered) from B1 onto the newly
# added header cells C1:F1.
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (rows 2-11): deathPlace / parent / deathDate populated,
#     birthDate / birthPlace left blank ---
$ws.Range("A2").Value = "http://dbpedia.org/resource/Giovanni_Francesco_Guidi_di_Bagno"
$ws.Range("B2").Value = "http://dbpedia.org/resource/Rome"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Colonna_family"
$ws.Range("D2").Value = "http://dbpedia.org/resource/1641"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""

$ws.Range("A3").Value = "http://dbpedia.org/resource/Giovanni_Doria"
$ws.Range("B3").Value = "http://dbpedia.org/resource/Palermo"
$ws.Range("C3").Value = "http://dbpedia.org/resource/Giovanni_Andrea_Doria"
$ws.Range("D3").Value = "http://dbpedia.org/resource/1642"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

$ws.Range("A4").Value = "http://dbpedia.org/resource/George_Murray_(bishop_of_Rochester)"
$ws.Range("B4").Value = "http://dbpedia.org/resource/Chester_Square"
$ws.Range("C4").Value = "http://dbpedia.org/resource/Lord_George_Murray_(bishop)"
$ws.Range("D4").Value = "http://dbpedia.org/resource/1860"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

$ws.Range("A5").Value = "http://dbpedia.org/resource/Geoffrey_(archbishop_of_York)"
$ws.Range("B5").Value = "http://dbpedia.org/resource/Normandy"
$ws.Range("C5").Value = "http://dbpedia.org/resource/Henry_II_of_England"
$ws.Range("D5").Value = "http://dbpedia.org/resource/12-12-12"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

$ws.Range("A6").Value = "http://dbpedia.org/resource/Ferdinand_III_of_Castile"
$ws.Range("B6").Value = "http://dbpedia.org/resource/Crown_of_Castile"
$ws.Range("C6").Value = "http://dbpedia.org/resource/Alfonso_IX_of_León"
$ws.Range("D6").Value = "http://dbpedia.org/resource/1252"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

$ws.Range("A7").Value = "http://dbpedia.org/resource/Erik_Benzelius_the_younger"
$ws.Range("B7").Value = "http://dbpedia.org/resource/Linköping"
$ws.Range("C7").Value = "http://dbpedia.org/resource/Erik_Benzelius_the_Elder"
$ws.Range("D7").Value = "http://dbpedia.org/resource/1743"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

$ws.Range("A8").Value = "http://dbpedia.org/resource/Edward_the_Confessor"
$ws.Range("B8").Value = "http://dbpedia.org/resource/London"
$ws.Range("C8").Value = "http://dbpedia.org/resource/Æthelred_the_Unready"
$ws.Range("D8").Value = "http://dbpedia.org/resource/1066"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

$ws.Range("A9").Value = "http://dbpedia.org/resource/Edward_William_Grinfield"
$ws.Range("B9").Value = "http://dbpedia.org/resource/Brighton"
$ws.Range("C9").Value = "http://dbpedia.org/resource/Thomas_Grinfield"
$ws.Range("D9").Value = "http://dbpedia.org/resource/1864"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

$ws.Range("A10").Value = "http://dbpedia.org/resource/Edward_Francis_Wilson"
$ws.Range("B10").Value = "http://dbpedia.org/resource/Saltspring_Island"
$ws.Range("C10").Value = "http://dbpedia.org/resource/Daniel_Wilson_(bishop)"
$ws.Range("D10").Value = "http://dbpedia.org/resource/1915"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("A11").Value = "http://dbpedia.org/resource/Donald_Foster_Hudson"
$ws.Range("B11").Value = "http://dbpedia.org/resource/England"
$ws.Range("C11").Value = "http://dbpedia.org/resource/Father"
$ws.Range("D11").Value = "http://dbpedia.org/resource/2003"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# --- Data rows (rows 12-18): deathPlace / parent / deathDate left blank,
#     birthDate / birthPlace populated ---
$ws.Range("A12").Value = "http://dbpedia.org/resource/Dick_Sheppard_(priest)"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "http://dbpedia.org/resource/1880"
$ws.Range("F12").Value = "http://dbpedia.org/resource/Windsor"

$ws.Range("A13").Value = "http://dbpedia.org/resource/Claus_Westermann"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "http://dbpedia.org/resource/1909"
$ws.Range("F13").Value = "http://dbpedia.org/resource/Berlin"

$ws.Range("A14").Value = "http://dbpedia.org/resource/Charles_Januarius_Acton"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "http://dbpedia.org/resource/1803"
$ws.Range("F14").Value = "http://dbpedia.org/resource/Naples"

$ws.Range("A15").Value = "http://dbpedia.org/resource/Carlo_Barberini"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = "http://dbpedia.org/resource/1630"
$ws.Range("F15").Value = "http://dbpedia.org/resource/Rome"

$ws.Range("A16").Value = "http://dbpedia.org/resource/Cardinal_de_Bouillon"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = "http://dbpedia.org/resource/1643"
$ws.Range("F16").Value = "http://dbpedia.org/resource/France"

$ws.Range("A17").Value = "http://dbpedia.org/resource/Camillo_Francesco_Maria_Pamphili"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = "http://dbpedia.org/resource/1622"
$ws.Range("F17").Value = "http://dbpedia.org/resource/Naples"

$ws.Range("A18").Value = "http://dbpedia.org/resource/Benjamin_Hoadly"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = "http://dbpedia.org/resource/1676"
$ws.Range("F18").Value = "http://dbpedia.org/resource/Kent"
